# làm giao diện QLDM - Tài khoản, thêm BM35,QD935, Thêm trong yêu cầu nghiệp vụ
#
# Adds row 37 (STT 35 "Danh mục tài khoản") to the "Yêu cầu nghiệp vụ" sheet,
# mirroring the layout/styling of row 34 (STT 32) on both halves of the
# table (A:D = yêu cầu nghiệp vụ, I:M = bảng trách nhiệm nghiệp vụ).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from row 34 (A:D and I:M — leave E:H alone, they are
# untouched in row 37 just like in every other row) down onto row 37 before
# writing the new values, so the new row keeps the same styles (number
# format / alignment) as the rest of the table.
$ws.Range("A34:D34").Copy()
$ws.Range("A37:D37").PasteSpecial(-4122)
$ws.Range("I34:M34").Copy()
$ws.Range("I37:M37").PasteSpecial(-4122)

# Left half of the table: "Yêu cầu nghiệp vụ"
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = "Danh mục tài khoản"
$ws.Range("C37").Value = "BM35"
$ws.Range("D37").Value = "QD35"

# Right half of the table: "Bảng trách nhiệm nghiệp vụ"
$ws.Range("I37").Value = 35
$ws.Range("J37").Value = "Danh mục tài khoản"
$ws.Range("K37").Value = "Cung cấp tài khoản cho nhân viên"
$ws.Range("L37").Value = "Xuất danh mục tài khoản và lưu "

# Match the author's final view state: scrolled down/right a bit, with
# L37 (the last edited cell) selected.
$ws.Range("L37").Select()
